$wb = $excel.ActiveWorkbook

# ALC (sheet1)
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(17, 8).Value = 1182.5  # H17: 1077.1666 -> 1182.5
$ws.Cells.Item(17, 10).Value = 1182.5  # J17: 1077.1666 -> 1182.5
$ws.Cells.Item(17, 12).Value = 3547.5  # L17: 3231.4998 -> 3547.5
$ws.Cells.Item(17, 14).Value = -3883.5  # N17: -3567.4998 -> -3883.5
$ws.Cells.Item(94, 8).Value = 1593.125  # H94: 1882.8334 -> 1593.125
$ws.Cells.Item(94, 9).Value = 1490.8334  # I94: 1882.8334 -> 1490.8334
$ws.Cells.Item(94, 10).Value = 1900  # J94: 0 -> 1900
$ws.Cells.Item(94, 11).Value = 1490.8334  # K94: 1882.8334 -> 1490.8334
$ws.Cells.Item(94, 12).Value = 1900  # L94: 0 -> 1900
$ws.Cells.Item(94, 13).Value = -1039.8334  # M94: -1431.8334 -> -1039.8334
$ws.Cells.Item(94, 14).Value = -2802  # N94: None -> -2802
$ws.Cells.Item(97, 8).Value = 10500  # H97: 16000 -> 10500
$ws.Cells.Item(97, 10).Value = 10500  # J97: 16000 -> 10500
$ws.Cells.Item(97, 12).Value = 31500  # L97: 48000 -> 31500
$ws.Cells.Item(97, 14).Value = -32492  # N97: -48992 -> -32492
$ws.Cells.Item(113, 8).Value = 6000  # H113: 4384.9 -> 6000
$ws.Cells.Item(113, 9).Value = 0  # I113: 4026.6667 -> 0
$ws.Cells.Item(113, 10).Value = 6000  # J113: 4538.4287 -> 6000
$ws.Cells.Item(113, 11).Value = 0  # K113: 4026.6667 -> 0
$ws.Cells.Item(113, 12).Value = 6000  # L113: 4538.4287 -> 6000
$ws.Cells.Item(113, 13).ClearContents()  # M113: -772.6667000000002 -> (removed)
$ws.Cells.Item(113, 14).Value = -12508  # N113: -11046.4287 -> -12508
$ws.Cells.Item(129, 8).Value = 1323.5714  # H129: 1280.8334 -> 1323.5714
$ws.Cells.Item(129, 9).Value = 853  # I129: 683.75 -> 853
$ws.Cells.Item(129, 10).Value = 2500  # J129: 2475 -> 2500
$ws.Cells.Item(129, 11).Value = 2559  # K129: 2051.25 -> 2559
$ws.Cells.Item(129, 12).Value = 7500  # L129: 7425 -> 7500
$ws.Cells.Item(129, 13).Value = 2441  # M129: 2948.75 -> 2441
$ws.Cells.Item(129, 14).Value = -17500  # N129: -17425 -> -17500
$ws.Cells.Item(132, 8).Value = 3143.9473  # H132: 4191.3076 -> 3143.9473
$ws.Cells.Item(132, 9).Value = 3160.8823  # I132: 4407.909 -> 3160.8823
$ws.Cells.Item(132, 11).Value = 9482.6469  # K132: 13223.727 -> 9482.6469
$ws.Cells.Item(132, 13).Value = -6952.6469  # M132: -10693.727 -> -6952.6469
$ws.Cells.Item(135, 8).Value = 1512.862  # H135: 1399.6897 -> 1512.862
$ws.Cells.Item(135, 9).Value = 1565.591  # I135: 1357.6666 -> 1565.591
$ws.Cells.Item(135, 10).Value = 1347.1428  # J135: 1510 -> 1347.1428
$ws.Cells.Item(135, 11).Value = 14090.319  # K135: 12218.9994 -> 14090.319
$ws.Cells.Item(135, 12).Value = 12124.2852  # L135: 13590 -> 12124.2852
$ws.Cells.Item(135, 13).Value = -11555.319  # M135: -9683.999400000001 -> -11555.319
$ws.Cells.Item(135, 14).Value = -17194.2852  # N135: -18660 -> -17194.2852
$ws.Cells.Item(138, 8).Value = 2941.1304  # H138: 3146.8572 -> 2941.1304
$ws.Cells.Item(138, 9).Value = 1079.75  # I138: 1248.2 -> 1079.75
$ws.Cells.Item(138, 11).Value = 3239.25  # K138: 3744.6 -> 3239.25
$ws.Cells.Item(138, 13).Value = 1900.75  # M138: 1395.4 -> 1900.75
$ws.Cells.Item(141, 8).Value = 4384.6  # H141: 5844.56 -> 4384.6
$ws.Cells.Item(141, 9).Value = 4176.3335  # I141: 3008.8462 -> 4176.3335
$ws.Cells.Item(141, 10).Value = 4576.846  # J141: 8916.583000000001 -> 4576.846
$ws.Cells.Item(141, 11).Value = 12529.0005  # K141: 9026.5386 -> 12529.0005
$ws.Cells.Item(141, 12).Value = 13730.538  # L141: 26749.749 -> 13730.538
$ws.Cells.Item(141, 13).Value = -7349.000499999998  # M141: -3846.5386 -> -7349.000499999998
$ws.Cells.Item(141, 14).Value = -24090.538  # N141: -37109.749 -> -24090.538

# ARM (sheet2)
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(74, 8).Value = 3079.1  # H74: 2569.0557 -> 3079.1
$ws.Cells.Item(74, 9).Value = 3079.1  # I74: 2569.0557 -> 3079.1
$ws.Cells.Item(74, 11).Value = 3079.1  # K74: 2569.0557 -> 3079.1
$ws.Cells.Item(74, 13).Value = -2205.1  # M74: -1695.0557 -> -2205.1
$ws.Cells.Item(77, 8).Value = 3079.1  # H77: 2569.0557 -> 3079.1
$ws.Cells.Item(77, 9).Value = 3079.1  # I77: 2569.0557 -> 3079.1
$ws.Cells.Item(77, 11).Value = 15395.5  # K77: 12845.2785 -> 15395.5
$ws.Cells.Item(77, 13).Value = -11027.5  # M77: -8477.2785 -> -11027.5
$ws.Cells.Item(122, 8).Value = 2237.0356  # H122: 2845.875 -> 2237.0356
$ws.Cells.Item(122, 9).Value = 2069.611  # I122: 2720.4285 -> 2069.611
$ws.Cells.Item(122, 10).Value = 2538.4  # J122: 2943.4443 -> 2538.4
$ws.Cells.Item(122, 11).Value = 6208.833  # K122: 8161.2855 -> 6208.833
$ws.Cells.Item(122, 12).Value = 7615.200000000001  # L122: 8830.332900000001 -> 7615.200000000001
$ws.Cells.Item(122, 13).Value = -3758.833  # M122: -5711.2855 -> -3758.833
$ws.Cells.Item(122, 14).Value = -12515.2  # N122: -13730.3329 -> -12515.2
$ws.Cells.Item(132, 8).Value = 2090.8  # H132: 1928.4517 -> 2090.8
$ws.Cells.Item(132, 9).Value = 1584.3636  # I132: 1650.6207 -> 1584.3636
$ws.Cells.Item(132, 10).Value = 5804.6665  # J132: 5957 -> 5804.6665
$ws.Cells.Item(132, 11).Value = 4753.0908  # K132: 4951.8621 -> 4753.0908
$ws.Cells.Item(132, 12).Value = 17413.9995  # L132: 17871 -> 17413.9995
$ws.Cells.Item(132, 13).Value = -2223.0908  # M132: -2421.8621 -> -2223.0908
$ws.Cells.Item(132, 14).Value = -22473.9995  # N132: -22931 -> -22473.9995

# BSM (sheet3)
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(96, 8).Value = 18392.143  # H96: 25466.875 -> 18392.143
$ws.Cells.Item(96, 10).Value = 75000  # J96: 74995 -> 75000
$ws.Cells.Item(96, 12).Value = 75000  # L96: 74995 -> 75000
$ws.Cells.Item(96, 14).Value = -80492  # N96: -80487 -> -80492
$ws.Cells.Item(105, 8).Value = 3431.6667  # H105: 1784.931 -> 3431.6667
$ws.Cells.Item(105, 9).Value = 2736.6667  # I105: 1153.0476 -> 2736.6667
$ws.Cells.Item(105, 10).Value = 3779.1667  # J105: 3443.625 -> 3779.1667
$ws.Cells.Item(105, 11).Value = 2736.6667  # K105: 1153.0476 -> 2736.6667
$ws.Cells.Item(105, 12).Value = 3779.1667  # L105: 3443.625 -> 3779.1667
$ws.Cells.Item(105, 13).Value = -989.6667000000002  # M105: 593.9523999999999 -> -989.6667000000002
$ws.Cells.Item(105, 14).Value = -7273.1667  # N105: -6937.625 -> -7273.1667
$ws.Cells.Item(134, 8).Value = 5595.6665  # H134: 3110.2273 -> 5595.6665
$ws.Cells.Item(134, 9).Value = 5595.6665  # I134: 3112.5 -> 5595.6665
$ws.Cells.Item(134, 10).Value = 0  # J134: 3100 -> 0
$ws.Cells.Item(134, 11).Value = 16786.9995  # K134: 9337.5 -> 16786.9995
$ws.Cells.Item(134, 12).Value = 0  # L134: 9300 -> 0
$ws.Cells.Item(134, 13).Value = -14251.9995  # M134: -6802.5 -> -14251.9995
$ws.Cells.Item(134, 14).ClearContents()  # N134: -14370 -> (removed)

# CRP (sheet4)
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 3110.3384  # H31: 3313.8196 -> 3110.3384
$ws.Cells.Item(31, 9).Value = 1885.7693  # I31: 2161 -> 1885.7693
$ws.Cells.Item(31, 10).Value = 3926.718  # J31: 4011.5789 -> 3926.718
$ws.Cells.Item(31, 11).Value = 1885.7693  # K31: 2161 -> 1885.7693
$ws.Cells.Item(31, 12).Value = 3926.718  # L31: 4011.5789 -> 3926.718
$ws.Cells.Item(31, 13).Value = -1590.7693  # M31: -1866 -> -1590.7693
$ws.Cells.Item(31, 14).Value = -4516.718  # N31: -4601.5789 -> -4516.718
$ws.Cells.Item(34, 8).Value = 3110.3384  # H34: 3313.8196 -> 3110.3384
$ws.Cells.Item(34, 9).Value = 1885.7693  # I34: 2161 -> 1885.7693
$ws.Cells.Item(34, 10).Value = 3926.718  # J34: 4011.5789 -> 3926.718
$ws.Cells.Item(34, 11).Value = 1885.7693  # K34: 2161 -> 1885.7693
$ws.Cells.Item(34, 12).Value = 3926.718  # L34: 4011.5789 -> 3926.718
$ws.Cells.Item(34, 13).Value = -1683.7693  # M34: -1959 -> -1683.7693
$ws.Cells.Item(34, 14).Value = -4330.718  # N34: -4415.5789 -> -4330.718
$ws.Cells.Item(99, 8).Value = 2186.5334  # H99: 2710.7144 -> 2186.5334
$ws.Cells.Item(99, 9).Value = 2210.6  # I99: 2649 -> 2210.6
$ws.Cells.Item(99, 10).Value = 2138.4  # J99: 2865 -> 2138.4
$ws.Cells.Item(99, 11).Value = 2210.6  # K99: 2649 -> 2210.6
$ws.Cells.Item(99, 12).Value = 2138.4  # L99: 2865 -> 2138.4
$ws.Cells.Item(99, 13).Value = -712.5999999999999  # M99: -1151 -> -712.5999999999999
$ws.Cells.Item(99, 14).Value = -5134.4  # N99: -5861 -> -5134.4
$ws.Cells.Item(122, 8).Value = 1885.75  # H122: 1883.7 -> 1885.75
$ws.Cells.Item(122, 9).Value = 1446.1111  # I122: 1357.3334 -> 1446.1111
$ws.Cells.Item(122, 10).Value = 3204.6667  # J122: 3462.8 -> 3204.6667
$ws.Cells.Item(122, 11).Value = 4338.3333  # K122: 4072.0002 -> 4338.3333
$ws.Cells.Item(122, 12).Value = 9614.000100000001  # L122: 10388.4 -> 9614.000100000001
$ws.Cells.Item(122, 13).Value = -1888.3333  # M122: -1622.0002 -> -1888.3333
$ws.Cells.Item(122, 14).Value = -14514.0001  # N122: -15288.4 -> -14514.0001
$ws.Cells.Item(126, 8).Value = 2186.5334  # H126: 2710.7144 -> 2186.5334
$ws.Cells.Item(126, 9).Value = 2210.6  # I126: 2649 -> 2210.6
$ws.Cells.Item(126, 10).Value = 2138.4  # J126: 2865 -> 2138.4
$ws.Cells.Item(126, 11).Value = 6631.799999999999  # K126: 7947 -> 6631.799999999999
$ws.Cells.Item(126, 12).Value = 6415.200000000001  # L126: 8595 -> 6415.200000000001
$ws.Cells.Item(126, 13).Value = -4161.799999999999  # M126: -5477 -> -4161.799999999999
$ws.Cells.Item(126, 14).Value = -11355.2  # N126: -13535 -> -11355.2
$ws.Cells.Item(127, 8).Value = 0  # H127: 50000 -> 0
$ws.Cells.Item(127, 9).Value = 0  # I127: 50000 -> 0
$ws.Cells.Item(127, 11).Value = 0  # K127: 50000 -> 0
$ws.Cells.Item(127, 13).ClearContents()  # M127: -45040 -> (removed)
$ws.Cells.Item(132, 8).Value = 2666.5  # H132: 1603.2727 -> 2666.5
$ws.Cells.Item(132, 9).Value = 2799.8  # I132: 1603.2727 -> 2799.8
$ws.Cells.Item(132, 10).Value = 2000  # J132: 0 -> 2000
$ws.Cells.Item(132, 11).Value = 8399.400000000001  # K132: 4809.8181 -> 8399.400000000001
$ws.Cells.Item(132, 12).Value = 6000  # L132: 0 -> 6000
$ws.Cells.Item(132, 13).Value = -5869.400000000001  # M132: -2279.8181 -> -5869.400000000001
$ws.Cells.Item(132, 14).Value = -11060  # N132: None -> -11060
$ws.Cells.Item(134, 8).Value = 2858.5833  # H134: 3470.5557 -> 2858.5833
$ws.Cells.Item(134, 9).Value = 3085.3  # I134: 3685.625 -> 3085.3
$ws.Cells.Item(134, 10).Value = 1725  # J134: 1750 -> 1725
$ws.Cells.Item(134, 11).Value = 9255.900000000001  # K134: 11056.875 -> 9255.900000000001
$ws.Cells.Item(134, 12).Value = 5175  # L134: 5250 -> 5175
$ws.Cells.Item(134, 13).Value = -6720.900000000001  # M134: -8521.875 -> -6720.900000000001
$ws.Cells.Item(134, 14).Value = -10245  # N134: -10320 -> -10245

# CUL (sheet5)
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(68, 8).Value = 4844.5454  # H68: 4769.9165 -> 4844.5454
$ws.Cells.Item(68, 10).Value = 5750  # J68: 5492.7144 -> 5750
$ws.Cells.Item(68, 12).Value = 17250  # L68: 16478.1432 -> 17250
$ws.Cells.Item(68, 14).Value = -18872  # N68: -18100.1432 -> -18872
$ws.Cells.Item(71, 8).Value = 4844.5454  # H71: 4769.9165 -> 4844.5454
$ws.Cells.Item(71, 10).Value = 5750  # J71: 5492.7144 -> 5750
$ws.Cells.Item(71, 12).Value = 51750  # L71: 49434.4296 -> 51750
$ws.Cells.Item(71, 14).Value = -59862  # N71: -57546.4296 -> -59862
$ws.Cells.Item(107, 8).Value = 2156.2666  # H107: 2250.4285 -> 2156.2666
$ws.Cells.Item(107, 9).Value = 4450.2  # I107: 3371.4285 -> 4450.2
$ws.Cells.Item(107, 10).Value = 1009.3  # J107: 1129.4286 -> 1009.3
$ws.Cells.Item(107, 11).Value = 13350.6  # K107: 10114.2855 -> 13350.6
$ws.Cells.Item(107, 12).Value = 3027.9  # L107: 3388.2858 -> 3027.9
$ws.Cells.Item(107, 13).Value = -11430.6  # M107: -8194.2855 -> -11430.6
$ws.Cells.Item(107, 14).Value = -6867.9  # N107: -7228.2858 -> -6867.9
$ws.Cells.Item(112, 8).Value = 201017.8  # H112: 127473.75 -> 201017.8
$ws.Cells.Item(112, 10).Value = 0  # J112: 4900.3335 -> 0
$ws.Cells.Item(112, 12).Value = 0  # L112: 14701.0005 -> 0
$ws.Cells.Item(112, 14).ClearContents()  # N112: -16917.0005 -> (removed)
$ws.Cells.Item(118, 8).Value = 1708.8  # H118: 2232 -> 1708.8
$ws.Cells.Item(118, 9).Value = 1136  # I118: 1464.5 -> 1136
$ws.Cells.Item(118, 10).Value = 4000  # J118: 2999.5 -> 4000
$ws.Cells.Item(118, 11).Value = 3408  # K118: 4393.5 -> 3408
$ws.Cells.Item(118, 12).Value = 12000  # L118: 8998.5 -> 12000
$ws.Cells.Item(118, 13).Value = -2165  # M118: -3150.5 -> -2165
$ws.Cells.Item(118, 14).Value = -14486  # N118: -11484.5 -> -14486
$ws.Cells.Item(121, 8).Value = 100785.3  # H121: 702 -> 100785.3
$ws.Cells.Item(121, 9).Value = 830.8  # I121: 586.8182 -> 830.8
$ws.Cells.Item(121, 10).Value = 200739.8  # J121: 860.375 -> 200739.8
$ws.Cells.Item(121, 11).Value = 2492.4  # K121: 1760.4546 -> 2492.4
$ws.Cells.Item(121, 12).Value = 602219.3999999999  # L121: 2581.125 -> 602219.3999999999
$ws.Cells.Item(121, 13).Value = -1182.4  # M121: -450.4546 -> -1182.4
$ws.Cells.Item(121, 14).Value = -604839.3999999999  # N121: -5201.125 -> -604839.3999999999
$ws.Cells.Item(136, 8).Value = 1200  # H136: 678.4 -> 1200
$ws.Cells.Item(136, 9).Value = 1200  # I136: 678.4 -> 1200
$ws.Cells.Item(136, 11).Value = 3600  # K136: 2035.2 -> 3600
$ws.Cells.Item(136, 13).Value = 1500  # M136: 3064.8 -> 1500

# GSM (sheet6)
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(48, 8).Value = 26999.5  # H48: 25999.666 -> 26999.5
$ws.Cells.Item(48, 10).Value = 26999.5  # J48: 25999.666 -> 26999.5
$ws.Cells.Item(48, 12).Value = 26999.5  # L48: 25999.666 -> 26999.5
$ws.Cells.Item(48, 14).Value = -27969.5  # N48: -26969.666 -> -27969.5
$ws.Cells.Item(122, 8).Value = 11632  # H122: 9639.477000000001 -> 11632
$ws.Cells.Item(122, 9).Value = 12619.125  # I122: 9639.477000000001 -> 12619.125
$ws.Cells.Item(122, 10).Value = 8999.666999999999  # J122: 0 -> 8999.666999999999
$ws.Cells.Item(122, 11).Value = 37857.375  # K122: 28918.431 -> 37857.375
$ws.Cells.Item(122, 12).Value = 26999.001  # L122: 0 -> 26999.001
$ws.Cells.Item(122, 13).Value = -35407.375  # M122: -26468.431 -> -35407.375
$ws.Cells.Item(122, 14).Value = -31899.001  # N122: None -> -31899.001
$ws.Cells.Item(132, 8).Value = 4366.6665  # H132: 2526.25 -> 4366.6665
$ws.Cells.Item(132, 9).Value = 4366.6665  # I132: 2360.6296 -> 4366.6665
$ws.Cells.Item(132, 10).Value = 0  # J132: 6998 -> 0
$ws.Cells.Item(132, 11).Value = 13099.9995  # K132: 7081.888800000001 -> 13099.9995
$ws.Cells.Item(132, 12).Value = 0  # L132: 20994 -> 0
$ws.Cells.Item(132, 13).Value = -10569.9995  # M132: -4551.888800000001 -> -10569.9995
$ws.Cells.Item(132, 14).ClearContents()  # N132: -26054 -> (removed)

# LTW (sheet7)
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(19, 8).Value = 0  # H19: 10200 -> 0
$ws.Cells.Item(19, 10).Value = 0  # J19: 10200 -> 0
$ws.Cells.Item(19, 12).Value = 0  # L19: 10200 -> 0
$ws.Cells.Item(19, 14).ClearContents()  # N19: -10540 -> (removed)
$ws.Cells.Item(82, 8).Value = 2823.92  # H82: 2546.077 -> 2823.92
$ws.Cells.Item(82, 9).Value = 1980.8572  # I82: 1834.6957 -> 1980.8572
$ws.Cells.Item(82, 10).Value = 7250  # J82: 8000 -> 7250
$ws.Cells.Item(82, 11).Value = 1980.8572  # K82: 1834.6957 -> 1980.8572
$ws.Cells.Item(82, 12).Value = 7250  # L82: 8000 -> 7250
$ws.Cells.Item(82, 13).Value = -1619.8572  # M82: -1473.6957 -> -1619.8572
$ws.Cells.Item(82, 14).Value = -7972  # N82: -8722 -> -7972
$ws.Cells.Item(85, 8).Value = 2823.92  # H85: 2546.077 -> 2823.92
$ws.Cells.Item(85, 9).Value = 1980.8572  # I85: 1834.6957 -> 1980.8572
$ws.Cells.Item(85, 10).Value = 7250  # J85: 8000 -> 7250
$ws.Cells.Item(85, 11).Value = 1980.8572  # K85: 1834.6957 -> 1980.8572
$ws.Cells.Item(85, 12).Value = 7250  # L85: 8000 -> 7250
$ws.Cells.Item(85, 13).Value = -732.8571999999999  # M85: -586.6957 -> -732.8571999999999
$ws.Cells.Item(85, 14).Value = -9746  # N85: -10496 -> -9746
$ws.Cells.Item(96, 8).Value = 0  # H96: 40000 -> 0
$ws.Cells.Item(96, 10).Value = 0  # J96: 40000 -> 0
$ws.Cells.Item(96, 12).Value = 0  # L96: 40000 -> 0
$ws.Cells.Item(96, 14).ClearContents()  # N96: -45492 -> (removed)

# WVR (sheet8)
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(119, 8).Value = 79995  # H119: 95000 -> 79995
$ws.Cells.Item(119, 10).Value = 79995  # J119: 95000 -> 79995
$ws.Cells.Item(119, 12).Value = 79995  # L119: 95000 -> 79995
$ws.Cells.Item(119, 14).Value = -89671  # N119: -104676 -> -89671
$ws.Cells.Item(122, 8).Value = 4394.591  # H122: 4929.45 -> 4394.591
$ws.Cells.Item(122, 9).Value = 1656.3334  # I122: 1999.1111 -> 1656.3334
$ws.Cells.Item(122, 10).Value = 7680.5  # J122: 7327 -> 7680.5
$ws.Cells.Item(122, 11).Value = 4969.0002  # K122: 5997.3333 -> 4969.0002
$ws.Cells.Item(122, 12).Value = 23041.5  # L122: 21981 -> 23041.5
$ws.Cells.Item(122, 13).Value = -2519.0002  # M122: -3547.3333 -> -2519.0002
$ws.Cells.Item(122, 14).Value = -27941.5  # N122: -26881 -> -27941.5
$ws.Cells.Item(124, 8).Value = 57999.75  # H124: 62500 -> 57999.75
$ws.Cells.Item(124, 10).Value = 57999.75  # J124: 62500 -> 57999.75
$ws.Cells.Item(124, 12).Value = 57999.75  # L124: 62500 -> 57999.75
$ws.Cells.Item(124, 14).Value = -67819.75  # N124: -72320 -> -67819.75
$ws.Cells.Item(132, 8).Value = 4973.4287  # H132: 4853.8335 -> 4973.4287
$ws.Cells.Item(132, 9).Value = 4973.4287  # I132: 5007.6787 -> 4973.4287
$ws.Cells.Item(132, 10).Value = 0  # J132: 2700 -> 0
$ws.Cells.Item(132, 11).Value = 14920.2861  # K132: 15023.0361 -> 14920.2861
$ws.Cells.Item(132, 12).Value = 0  # L132: 8100 -> 0
$ws.Cells.Item(132, 13).Value = -12390.2861  # M132: -12493.0361 -> -12390.2861
$ws.Cells.Item(132, 14).ClearContents()  # N132: -13160 -> (removed)
$ws.Cells.Item(141, 8).Value = 25000  # H141: 100000 -> 25000
$ws.Cells.Item(141, 10).Value = 25000  # J141: 100000 -> 25000
$ws.Cells.Item(141, 12).Value = 25000  # L141: 100000 -> 25000
$ws.Cells.Item(141, 14).Value = -35360  # N141: -110360 -> -35360
